$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 9 with the new variable entry (3rd plot)
$ws.Range("A9").Value = "PreProcessDone"
$ws.Range("B9").Value = "True False"
$ws.Range("C9").Value = "data has been pre processed"

# Update the active selection to C9 to match the saved view state
$ws.Range("C9").Select()
